# Add single family home data to the RICS service life table.
# Strategy: insert the 8 new rows (processed bottom-to-top by original
# anchor row so earlier row numbers stay valid), then (re)write every
# data cell's value explicitly, and finally fix up the handful of
# cells whose style must be "no direct format" (vs the inherited "s=1").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1. Insert the 8 new rows (bottom-to-top on ORIGINAL row numbers) ----

# New row after the old row 22 (EPDM roofing) -> will hold "Asphalt shingle roofing"
$ws.Rows.Item(23).Insert()

# New row before the old row 22 (EPDM roofing) -> will hold "Glazing: operable window"
$ws.Rows.Item(22).Insert()

# New row before the old row 14 (Insulated metal panel) -> "Brick: wood framing"
$ws.Rows.Item(14).Insert()

# 4 new rows before the old row 10 (Curtain wall: steel spandrel) ->
# "Floor framing", "Sub-flooring", "Roof framing", "Roof decking"
$ws.Rows.Item(10).Resize(4).Insert()

# New row before the old row 3 (Wall foundation) -> "Concrete footing"
$ws.Rows.Item(3).Insert()

# ---- 2. Write the full A:D content for every data row (2-30) ----

$data = @(
  @{Row=2;  A=1;  B='Column foundation';               C='piling and foundations';                          D=60},
  @{Row=3;  A=2;  B='Concrete footing';                 C='piling and foundations';                          D=60},
  @{Row=4;  A=3;  B='Wall foundation';                  C='piling and foundations';                          D=60},
  @{Row=5;  A=4;  B='Slab on grade';                    C='Lowest ground floor';                             D=60},
  @{Row=6;  A=5;  B='Elevated slabs';                   C='Superstructure: Structural elements';             D=60},
  @{Row=7;  A=6;  B='Structural framing: beams';        C='Superstructure: Structural elements';             D=60},
  @{Row=8;  A=7;  B='Structural framing: girders';      C='Superstructure: Structural elements';             D=60},
  @{Row=9;  A=8;  B='Structural columns';               C='Superstructure: Structural elements';             D=60},
  @{Row=10; A=9;  B='Structural walls';                 C='Superstructure: Structural elements';             D=60},
  @{Row=11; A=10; B='Floor framing';                    C='Superstructure: Structural elements';             D=60},
  @{Row=12; A=11; B='Sub-flooring';                     C='Superstructure: Structural elements';             D=60},
  @{Row=13; A=12; B='Roof framing';                     C='Superstructure: Structural elements';             D=60},
  @{Row=14; A=13; B='Roof decking';                     C='Superstructure: Structural elements';             D=60},
  @{Row=15; A=14; B='Curtain wall: steel spandrel';     C='Glazed cladded/curtain walling';                  D=35},
  @{Row=16; A=15; B='Curtain wall: aluminum spandrel';  C='Glazed cladded/curtain walling';                  D=35},
  @{Row=17; A=16; B='MV: brick';                        C='Brick, stone, block and precast concrete panels'; D=60},
  @{Row=18; A=17; B='MV: granite';                      C='Brick, stone, block and precast concrete panels'; D=60},
  @{Row=19; A=18; B='Brick: wood framing';              C='Brick, stone, block and precast concrete panels'; D=60},
  @{Row=20; A=19; B='Insulated metal panel';            C='Internal partitioning and dry lining: studwork';  D=30},
  @{Row=21; A=20; B='EIFS (XPS)';                       C='Glazed cladded/curtain walling';                  D=35},
  @{Row=22; A=21; B='Rainscreen: GFRC';                 C='Rain screens, timber panels';                     D=30},
  @{Row=23; A=22; B='Rainscreen: thin brick';           C='Rain screens, timber panels';                     D=30},
  @{Row=24; A=23; B='Rainscreen: wood';                 C='Rain screens, timber panels';                     D=30},
  @{Row=25; A=24; B='Rainscreen: formed steel panel';   C='Rain screens, timber panels';                     D=30},
  @{Row=26; A=25; B='Glazing: double pane IGU';         C='Glazed cladded/curtain walling';                  D=35},
  @{Row=27; A=26; B='Glazing: triple pane IGU';         C='Glazed cladded/curtain walling';                  D=35},
  @{Row=28; A=27; B='Glazing: operable window';         C='Glazed cladded/curtain walling';                  D=35},
  @{Row=29; A=28; B='EPDM roofing';                     C='Roof covering: Single-ply membrane';              D=30},
  @{Row=30; A=29; B='Asphalt shingle roofing';          C='Roof covering: Single-ply membrane';              D=30}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
}

# ---- 3. Fix up direct formatting on the cells that must carry no
#         explicit style (rows 3, 19, 28, 30 in the final layout) ----

$ws.Range("B3").Style = "Normal"
$ws.Range("B19").Style = "Normal"
$ws.Range("B28").Style = "Normal"
$ws.Range("B30").Style = "Normal"

# ---- 4. Selection + dimension bookkeeping (matches the authored diff) ----

$ws.Range("C30").Select()
